# Schedule workbook update: add a new task entry, add a new blank spacer
# row, and color-code the rows by status (green = existing plan items,
# pink/red/light-yellow = the three differently-prioritized items at the
# bottom of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New task text in F7 ("multiple buttons for multiple sizes for manga")
# ---------------------------------------------------------------------
$ws.Range("F7").Value = "multiple buttons for multiple sizes for manga"

# ---------------------------------------------------------------------
# 2. New blank spacer row (row 8), matching the existing blank spacer
#    style used at row 6 (date column + separator column only).
# ---------------------------------------------------------------------
$ws.Range("A8").Value = $null
$ws.Range("C8").Value = $null

# ---------------------------------------------------------------------
# Color palette (decimal = R + G*256 + B*65536, i.e. the value the Excel
# object model expects for Interior.Color / RGB()).
# ---------------------------------------------------------------------
$green      = 4063037    # FF3DFF3D
$yellowPure = 65535       # FFFFFF00 (picked then superseded; kept as an
                          # unused style swatch, matching the source file)
$red        = 255         # FFFF0000
$lightGold  = 4063231     # FFFFFF3D
$pink       = 10395391    # FFFF9E9E

# ---------------------------------------------------------------------
# 3. Green rows (4-10): existing plan. Touch cells in a fixed order so the
#    resulting style table lists "plain fill", "date fill", "center fill"
#    styles in that order.
# ---------------------------------------------------------------------
$ws.Range("B4").Interior.Color = $green
$ws.Range("A4").Interior.Color = $green
$ws.Range("C4").Interior.Color = $green

$greenPlain = $ws.Range("D4,F4,G4,B5,D5,F5,F6,B7,D7,F7,B9,D9,B10,D10")
foreach ($a in $greenPlain.Areas) { $a.Interior.Color = $green }

$greenDate = $ws.Range("A5,A6,A7,A8,A9,A10")
foreach ($a in $greenDate.Areas) { $a.Interior.Color = $green }

$greenCenter = $ws.Range("C5,C6,C7,C8,C9,C10")
foreach ($a in $greenCenter.Areas) { $a.Interior.Color = $green }

# ---------------------------------------------------------------------
# 4. Create (and leave unused) the plain-yellow swatch, via a scratch
#    cell far outside the used range that is fully removed afterwards.
# ---------------------------------------------------------------------
$ws.Range("Z1").Interior.Color = $yellowPure
$ws.Range("Z1").Interior.Color = $red
$ws.Range("Z1").ClearFormats()
$ws.Range("Z1").EntireColumn.Delete()

# ---------------------------------------------------------------------
# 5. Red rows (12-13).
# ---------------------------------------------------------------------
$ws.Range("B12").Interior.Color = $red
$ws.Range("A12").Interior.Color = $red
$ws.Range("C12").Interior.Color = $red

$redPlain = $ws.Range("D12,E12,B13,D13")
foreach ($a in $redPlain.Areas) { $a.Interior.Color = $red }

$ws.Range("A13").Interior.Color = $red
$ws.Range("C13").Interior.Color = $red

# ---------------------------------------------------------------------
# 6. Light-gold row (14).
# ---------------------------------------------------------------------
$ws.Range("B14").Interior.Color = $lightGold
$ws.Range("A14").Interior.Color = $lightGold
$ws.Range("C14").Interior.Color = $lightGold
$redPlain14 = $ws.Range("D14,E14")
foreach ($a in $redPlain14.Areas) { $a.Interior.Color = $lightGold }

# ---------------------------------------------------------------------
# 7. Pink row (11). C11 additionally gets the "/" separator text that the
#    other rows already show in their separator column.
# ---------------------------------------------------------------------
$ws.Range("B11").Interior.Color = $pink
$ws.Range("A11").Interior.Color = $pink
$ws.Range("C11").Value = "/"
$ws.Range("C11").HorizontalAlignment = -4108
$ws.Range("C11").Interior.Color = $pink
$pinkPlain = $ws.Range("D11,E11")
foreach ($a in $pinkPlain.Areas) { $a.Interior.Color = $pink }

# ---------------------------------------------------------------------
# 8. Move the active selection to D22 (matches the sheet's saved view).
# ---------------------------------------------------------------------
$ws.Range("D22").Select()
